$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-110 down to 52-111
$ws.Rows("51:51").Insert()

# Fill the new row 51 with data. Columns A,B,C,E,F,G,H,I,J,L,Q,R carry over the
# same values that were already present in the (now shifted) row below it.
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44658
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100109
$ws.Cells.Item(51, 8).Value = "Uva"
$ws.Cells.Item(51, 9).Value = 100109001
$ws.Cells.Item(51, 10).Value = "Uva"
$ws.Cells.Item(51, 11).Value = "Red Globe"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 300
$ws.Cells.Item(51, 14).Value = 9500
$ws.Cells.Item(51, 15).Value = 10000
$ws.Cells.Item(51, 16).Value = 9750
$ws.Cells.Item(51, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(51, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 19).Value = 542
$ws.Cells.Item(51, 20).Value = 18
